# Update "想去人数" (F column) figures for 展览 (sheet1) and 全部类型 (sheet4) worksheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")


# 展览 (sheet1)
$ws1.Range("F2").Value = 5538  # was 5529
$ws1.Range("F4").Value = 7651  # was 7640
$ws1.Range("F5").Value = 50  # was 49
$ws1.Range("F8").Value = 612  # was 611
$ws1.Range("F9").Value = 70  # was 68
$ws1.Range("F11").Value = 4421  # was 4402
$ws1.Range("F12").Value = 1792  # was 1787
$ws1.Range("F15").Value = 2986  # was 2979
$ws1.Range("F17").Value = 572  # was 571
$ws1.Range("F18").Value = 222  # was 221
$ws1.Range("F19").Value = 543  # was 540
$ws1.Range("F20").Value = 472  # was 470
$ws1.Range("F21").Value = 480  # was 478
$ws1.Range("F22").Value = 342  # was 338
$ws1.Range("F23").Value = 118  # was 116
$ws1.Range("F24").Value = 1722  # was 1719
$ws1.Range("F25").Value = 1249  # was 1244
$ws1.Range("F26").Value = 102  # was 101
$ws1.Range("F27").Value = 1443  # was 1437
$ws1.Range("F30").Value = 35  # was 34
$ws1.Range("F31").Value = 519  # was 518
$ws1.Range("F36").Value = 73  # was 72
$ws1.Range("F37").Value = 3103  # was 3081
$ws1.Range("F39").Value = 45  # was 44
$ws1.Range("F40").Value = 144  # was 139
$ws1.Range("F42").Value = 811  # was 797

# 全部类型 (sheet4)
$ws4.Range("F2").Value = 5538  # was 5529
$ws4.Range("F4").Value = 7651  # was 7640
$ws4.Range("F5").Value = 50  # was 49
$ws4.Range("F8").Value = 612  # was 611
$ws4.Range("F9").Value = 70  # was 68
$ws4.Range("F11").Value = 4421  # was 4402
$ws4.Range("F12").Value = 1792  # was 1787
$ws4.Range("F15").Value = 2986  # was 2979
$ws4.Range("F17").Value = 572  # was 571
$ws4.Range("F18").Value = 222  # was 221
$ws4.Range("F19").Value = 543  # was 540
$ws4.Range("F20").Value = 472  # was 470
$ws4.Range("F21").Value = 480  # was 478
$ws4.Range("F23").Value = 342  # was 338
$ws4.Range("F24").Value = 118  # was 116
$ws4.Range("F25").Value = 1722  # was 1719
$ws4.Range("F26").Value = 1249  # was 1244
$ws4.Range("F27").Value = 102  # was 101
$ws4.Range("F28").Value = 1443  # was 1437
$ws4.Range("F31").Value = 35  # was 34
$ws4.Range("F32").Value = 519  # was 518
$ws4.Range("F37").Value = 73  # was 72
$ws4.Range("F38").Value = 3103  # was 3082
$ws4.Range("F41").Value = 45  # was 44
$ws4.Range("F42").Value = 144  # was 139
$ws4.Range("F44").Value = 811  # was 797
